# Update the cryptos list with the latest scraped prices/volumes.
# GitHub Actions scheduled refresh.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($range, [string]$text)
    # Force the cell to stay plain text (NumberFormat "@") before assigning,
    # otherwise Excel auto-coerces numeric-looking strings like "243.17" or
    # "0.0000250" into real numbers and mangles/round-trips them (loses
    # trailing zeros, drops thousand-separator dots, etc). Reset the style
    # back to Normal afterwards so we don't leave a stray "@" number format
    # applied to the cell (these cells carry no style in the source sheet).
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.Style = "Normal"
}

# --- Rows whose coin/link stayed the same row, only Price (D) and/or
#     Volume(1h) (E) text changed -------------------------------------------------
$updates = @(
    @{ Row = 2;  D = "90.599.15";  E = "  +0.47%  " },
    @{ Row = 3;  D = "3.107.64";   E = "  +0.97%  " },
    @{ Row = 4;  D = $null;        E = "  +0.04%  " },
    @{ Row = 5;  D = "243.17";     E = "  +10.92%  " },
    @{ Row = 6;  D = "626.32";     E = "  +2.15%  " },
    @{ Row = 7;  D = $null;        E = "  +5.96%  " },
    @{ Row = 8;  D = "0.371";      E = "  +5.86%  " },
    @{ Row = 9;  D = $null;        E = "  +0.05%  " },
    @{ Row = 10; D = "3.108.14";   E = "  +1.12%  " },
    @{ Row = 11; D = $null;        E = "  +3.32%  " },
    @{ Row = 12; D = "0.203";      E = "  +4.16%  " },
    @{ Row = 13; D = "0.0000250";  E = "  +4.32%  " },
    @{ Row = 14; D = "35.50";      E = "  +3.48%  " },
    @{ Row = 15; D = "5.49";       E = "  -0.10%  " },
    @{ Row = 16; D = "90.477.90";  E = "  +0.58%  " },
    @{ Row = 17; D = "3.677.08";   E = "  +1.02%  " },

    @{ Row = 20; D = "14.32";      E = "  +1.19%  " },
    @{ Row = 21; D = "0.0000213";  E = "  +4.34%  " },
    @{ Row = 22; D = "5.75";       E = "  +7.75%  " },
    @{ Row = 23; D = "447.31";     E = "  +0.42%  " },
    @{ Row = 24; D = "9.02";       E = "  +2.02%  " },
    @{ Row = 25; D = "5.99";       E = "  +0.67%  " },
    @{ Row = 26; D = "92.61";      E = "  +2.18%  " },
    @{ Row = 27; D = "12.08";      E = "  +1.58%  " },
    @{ Row = 28; D = $null;        E = "  +0.93%  " },
    @{ Row = 29; D = $null;        E = "  +0.11%  " },
    @{ Row = 30; D = "0.175";      E = "  +10.20%  " },

    @{ Row = 33; D = $null;        E = "  +7.09%  " },
    @{ Row = 34; D = $null;        E = "  +31.25%  " },

    @{ Row = 36; D = $null;        E = "  -5.01%  " },

    @{ Row = 39; D = "1.93";       E = "  +2.24%  " },
    @{ Row = 40; D = "496.87";     E = "  -0.13%  " },
    @{ Row = 41; D = "3.66";       E = "  +6.80%  " },
    @{ Row = 42; D = "1.31";       E = "  +1.75%  " },
    @{ Row = 43; D = "0.418";      E = "  -1.99%  " },
    @{ Row = 44; D = "22.12";      E = "  -0.33%  " },
    @{ Row = 45; D = $null;        E = "  -0.01%  " },
    @{ Row = 46; D = "159.48";     E = "  +7.67%  " },
    @{ Row = 47; D = "1.92";       E = "  -0.73%  " },
    @{ Row = 48; D = "0.687";      E = "  -0.51%  " },
    @{ Row = 49; D = "4.55";       E = "  -0.94%  " },
    @{ Row = 50; D = "44.74";      E = "  +0.67%  " },
    @{ Row = 51; D = "1.34";       E = "  +1.52%  " }
)

foreach ($u in $updates) {
    $r = $u.Row
    if ($null -ne $u.D) {
        Set-TextValue $ws.Cells.Item($r, 4) $u.D
    }
    Set-TextValue $ws.Cells.Item($r, 5) $u.E
}

# --- Rows whose coin re-ranked (whole B/C/D/E row contents swapped) ------------

# Rows 18 & 19 swap: SuiNetwork now ranks above WrappedEther.
$ws.Cells.Item(18, 2).Value = "SuiNetwork"
$ws.Cells.Item(18, 3).Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
Set-TextValue $ws.Cells.Item(18, 4) "3.86"
Set-TextValue $ws.Cells.Item(18, 5) "  +4.66%  "

$ws.Cells.Item(19, 2).Value = "WrappedEther"
$ws.Cells.Item(19, 3).Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
Set-TextValue $ws.Cells.Item(19, 4) "3.093.62"
Set-TextValue $ws.Cells.Item(19, 5) "  +0.80%  "

# Rows 31 & 32 swap: InternetComputer(DFINITY) now ranks above Stellar.
$ws.Cells.Item(31, 2).Value = "InternetComputer(DFINITY)"
$ws.Cells.Item(31, 3).Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
Set-TextValue $ws.Cells.Item(31, 4) "9.30"
Set-TextValue $ws.Cells.Item(31, 5) "  -0.17%  "

$ws.Cells.Item(32, 2).Value = "Stellar"
$ws.Cells.Item(32, 3).Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
Set-TextValue $ws.Cells.Item(32, 4) "0.214"
Set-TextValue $ws.Cells.Item(32, 5) "  +5.97%  "

# Rows 35, 37, 38 rotate: RenderToken moves up to 35, MantraDAO moves to 37,
# Kaspa moves to 38 (row 36, EthereumClassic, stays put and only gets its
# Volume(1h) updated above).
$ws.Cells.Item(35, 2).Value = "RenderToken"
$ws.Cells.Item(35, 3).Value = "https://coinranking.com/coin/vfo5XYwcV+rendertoken-render"
Set-TextValue $ws.Cells.Item(35, 4) "7.82"
Set-TextValue $ws.Cells.Item(35, 5) "  +14.71%  "

$ws.Cells.Item(37, 2).Value = "MantraDAO"
$ws.Cells.Item(37, 3).Value = "https://coinranking.com/coin/cTdD8lD-6+mantradao-om"
Set-TextValue $ws.Cells.Item(37, 4) "4.29"
Set-TextValue $ws.Cells.Item(37, 5) "  +42.68%  "

$ws.Cells.Item(38, 2).Value = "Kaspa"
$ws.Cells.Item(38, 3).Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
Set-TextValue $ws.Cells.Item(38, 4) "0.156"
Set-TextValue $ws.Cells.Item(38, 5) "  +5.21%  "
